# Inserts a new data row at row 203 (pushing the existing rows 203-308 down
# to 204-309) and populates the new row with a "Primera" quality record for
# Cebollín at Femacal de La Calera, dated 2021-10-20 (Excel serial 44489).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("203:203").Insert()

$ws.Range("A203").Value = 3
$ws.Range("B203").Value = "Femacal de La Calera"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = "2021-10-20"
$ws.Range("E203").Value = 5
$ws.Range("F203").Value = 100112037
$ws.Range("G203").Value = "Cebollín"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 180
$ws.Range("K203").Value = 3000
$ws.Range("L203").Value = 3000
$ws.Range("M203").Value = 3000
$ws.Range("N203").Value = "$/paquete 36 unidades"
$ws.Range("O203").Value = "Provincia de Quillota"
$ws.Range("P203").Value = 83
$ws.Range("Q203").Value = 36
$ws.Range("R203").Value = "Hortaliza"
